# Appends four new rows to the "Arbeitsnachweis" table, describing the
# work done on the Fahrplanauskunft input mask / index page design.

function Set-CellRuns($Cell, $Paragraphs) {
    $ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $xmlParts = @()
    $expectedCount = $Paragraphs.Count
    foreach ($para in $Paragraphs) {
        $runsXml = ""
        foreach ($run in $para) {
            $brk = ""
            if ($run.ContainsKey("break") -and $run.break) {
                $brk = "<w:lastRenderedPageBreak/>"
            }
            $text = $run.text
            $preserve = ""
            if ($text -ne $text.Trim() -or $text -eq "") {
                $preserve = ' xml:space="preserve"'
            }
            $escaped = $text -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;'
            $runsXml += "<w:r>$brk<w:t$preserve>$escaped</w:t></w:r>"
        }
        $xmlParts += "<w:p $ns>$runsXml</w:p>"
    }
    $xml = [string]::Join("", $xmlParts)
    $Cell.Range.InsertXML($xml)
    # Depending on the target cell, InsertXML sometimes leaves the
    # original (now-empty) paragraph in place ahead of the inserted
    # content instead of merging into it - trim any such leftovers.
    while ($Cell.Range.Paragraphs.Count -gt $expectedCount) {
        $Cell.Range.Paragraphs.Item(1).Range.Delete()
    }
}

function Add-ArbeitsnachweisRow($Table, $Aufgabe, $Mitarbeiter, $Beschreibung, $Zeit, $Datum) {
    $row = $Table.Rows.Add()
    Set-CellRuns $row.Cells.Item(1) $Aufgabe
    Set-CellRuns $row.Cells.Item(2) $Mitarbeiter
    Set-CellRuns $row.Cells.Item(3) $Beschreibung
    Set-CellRuns $row.Cells.Item(4) $Zeit
    Set-CellRuns $row.Cells.Item(5) $Datum
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row: "Eingabe für Fahrplanauskunft erstellt" ---------------------
$aufgabe = @()
$aufgabe += ,@(@{text="Eingabe für Fahrplanauskunft erstellt"})

$mitarbeiter = @()
$mitarbeiter += ,@(@{text="Falk"})

$beschreibung = @()
$beschreibung += ,@(@{text="Responsive Maske zur Eingabe der Fahrplansuchkriterien "}, @{text="begonnen"})

$zeit = @()
$zeit += ,@(@{text="4h"})

$datum = @()
$datum += ,@(@{text="06.02.2023"})

Add-ArbeitsnachweisRow $t $aufgabe $mitarbeiter $beschreibung $zeit $datum

# --- Row: "Anfänge Suchlogik, Vorbereitung Übergabefelder" ------------
$aufgabe = @()
$aufgabe += ,@(@{text="Anfänge Suchlogik, Vorbereitung Übergabefelder"})

$mitarbeiter = @()
$mitarbeiter += ,@(@{text="Falk"})

$beschreibung = @()
$beschreibung += ,@(@{text="Request-Klassen erstellt, "}, @{text="PanelBuilder erstellt"})

$zeit = @()
$zeit += ,@(@{text="2h"})

$datum = @()
$datum += ,@(@{text="07.02.2023"})

Add-ArbeitsnachweisRow $t $aufgabe $mitarbeiter $beschreibung $zeit $datum

# --- Row: "Responsives UI, JavaScript Form modifizieren" --------------
$aufgabe = @()
$aufgabe += ,@(@{text="Responsives UI, JavaScript Form modifizieren"})

$mitarbeiter = @()
$mitarbeiter += ,@(@{text="Falk"})

$beschreibung = @()
$beschreibung += ,@(@{text="Eingabemaske für die Fahrplanauskunft fortgesetzt"}, @{text=", "}, @{text="Vorbereitung Fahrplananzeige "})

$zeit = @()
$zeit += ,@(@{text="3.5h"})

$datum = @()
$datum += ,@(@{text="15.02.2023"})

Add-ArbeitsnachweisRow $t $aufgabe $mitarbeiter $beschreibung $zeit $datum

# --- Row: "Design der Index-Seite begonnen" ---------------------------
$aufgabe = @()
$aufgabe += ,@(@{text="Design der Index-Seite begonnen"})

$mitarbeiter = @()
$mitarbeiter += ,@(@{text="Hollmann"})

$beschreibung = @()
$beschreibung += ,@(@{text="Kacheln und Grunddesign auf der Index-Seite erstellt (erste Idee)"})
$beschreibung += ,@(@{text="Fehlerbehebung in Fahrplanauskunft"; break=$true}, @{text=", kleinere Codeanpassungen"})

$zeit = @()
$zeit += ,@(@{text="1"; break=$true}, @{text=".5"}, @{text="h"})

$datum = @()
$datum += ,@(@{text="17.02.2023"})

Add-ArbeitsnachweisRow $t $aufgabe $mitarbeiter $beschreibung $zeit $datum
